$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C15").Value = "  "
$ws.Range("C16").Select()
